$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
# Row 107
$wsALC.Range("H107").Value = 645.0952
$wsALC.Range("I107").Value = 550
$wsALC.Range("J107").Value = 799.625
$wsALC.Range("K107").Value = 550
$wsALC.Range("L107").Value = 799.625
$wsALC.Range("M107").Value = 1370
$wsALC.Range("N107").Value = -4639.625
# Row 132
$wsALC.Range("H132").Value = 2658
$wsALC.Range("I132").Value = 1688.2106
$wsALC.Range("K132").Value = 5064.6318
$wsALC.Range("M132").Value = -2534.6318
# Row 137
$wsALC.Range("H137").Value = 2142.0815
$wsALC.Range("I137").Value = 2306.5625
$wsALC.Range("J137").Value = 1832.4706
$wsALC.Range("K137").Value = 6919.6875
$wsALC.Range("L137").Value = 5497.4118
$wsALC.Range("M137").Value = -4369.6875
$wsALC.Range("N137").Value = -10597.4118

$wsARM = $wb.Worksheets.Item("ARM")
# Row 32
$wsARM.Range("H32").Value = 1072112.5
$wsARM.Range("I32").Value = 1285407.9
$wsARM.Range("J32").Value = 5635.6924
$wsARM.Range("K32").Value = 1285407.9
$wsARM.Range("L32").Value = 5635.6924
$wsARM.Range("M32").Value = -1285120.9
$wsARM.Range("N32").Value = -6209.6924
# Row 40
$wsARM.Range("H40").Value = 14600
$wsARM.Range("J40").Value = 14600
$wsARM.Range("L40").Value = 14600
$wsARM.Range("N40").Value = -14952
# Row 42
$wsARM.Range("H42").Value = 9325
$wsARM.Range("J42").Value = 9325
$wsARM.Range("L42").Value = 9325
$wsARM.Range("N42").Value = -10297
# Row 46
$wsARM.Range("H46").Value = 2173
$wsARM.Range("I46").Value = 846
$wsARM.Range("J46").Value = 3500
$wsARM.Range("K46").Value = 846
$wsARM.Range("L46").Value = 3500
$wsARM.Range("N46").Value = -4138
$wsARM.Range("M46").Value = -527
# Row 61
$wsARM.Range("H61").Value = 660002.25
$wsARM.Range("I61").Value = 650321.1
$wsARM.Range("J61").Value = 670974.25
$wsARM.Range("K61").Value = 650321.1
$wsARM.Range("L61").Value = 670974.25
$wsARM.Range("M61").Value = -650109.1
$wsARM.Range("N61").Value = -671398.25
# Row 74
$wsARM.Range("H74").Value = 264000.9
$wsARM.Range("I74").Value = 385994.5
$wsARM.Range("J74").Value = 65761.31
$wsARM.Range("K74").Value = 385994.5
$wsARM.Range("L74").Value = 65761.31
$wsARM.Range("M74").Value = -385120.5
$wsARM.Range("N74").Value = -67509.31
# Row 77
$wsARM.Range("H77").Value = 264000.9
$wsARM.Range("I77").Value = 385994.5
$wsARM.Range("J77").Value = 65761.31
$wsARM.Range("K77").Value = 1929972.5
$wsARM.Range("L77").Value = 328806.55
$wsARM.Range("M77").Value = -1925604.5
$wsARM.Range("N77").Value = -337542.55
# Row 122
$wsARM.Range("H122").Value = 4420.6
$wsARM.Range("I122").Value = 3760.5
$wsARM.Range("J122").Value = 7061
$wsARM.Range("K122").Value = 11281.5
$wsARM.Range("L122").Value = 21183
$wsARM.Range("M122").Value = -8831.5
$wsARM.Range("N122").Value = -26083
# Row 132
$wsARM.Range("H132").Value = 27255.88
$wsARM.Range("I132").Value = 46979.824
$wsARM.Range("J132").Value = 3379.5264
$wsARM.Range("K132").Value = 140939.472
$wsARM.Range("L132").Value = 10138.5792
$wsARM.Range("M132").Value = -138409.472
$wsARM.Range("N132").Value = -15198.5792
# Row 136
$wsARM.Range("H136").Value = 660002.25
$wsARM.Range("I136").Value = 650321.1
$wsARM.Range("J136").Value = 670974.25
$wsARM.Range("K136").Value = 1950963.3
$wsARM.Range("L136").Value = 2012922.75
$wsARM.Range("M136").Value = -1948413.3
$wsARM.Range("N136").Value = -2018022.75

$wsBSM = $wb.Worksheets.Item("BSM")
# Row 5
$wsBSM.Range("H5").Value = 382
$wsBSM.Range("I5").Value = 382
$wsBSM.Range("K5").Value = 382
$wsBSM.Range("M5").Value = -269
# Row 99
$wsBSM.Range("H99").Value = 3335.7273
$wsBSM.Range("I99").Value = 4619.9287
$wsBSM.Range("J99").Value = 1088.375
$wsBSM.Range("K99").Value = 4619.9287
$wsBSM.Range("L99").Value = 1088.375
$wsBSM.Range("M99").Value = -3121.9287
$wsBSM.Range("N99").Value = -4084.375

$wsCRP = $wb.Worksheets.Item("CRP")
# Row 31
$wsCRP.Range("H31").Value = 2657.1875
$wsCRP.Range("I31").Value = 1955.0667
$wsCRP.Range("J31").Value = 3827.389
$wsCRP.Range("K31").Value = 1955.0667
$wsCRP.Range("L31").Value = 3827.389
$wsCRP.Range("M31").Value = -1660.0667
$wsCRP.Range("N31").Value = -4417.389
# Row 34
$wsCRP.Range("H34").Value = 2657.1875
$wsCRP.Range("I34").Value = 1955.0667
$wsCRP.Range("J34").Value = 3827.389
$wsCRP.Range("K34").Value = 1955.0667
$wsCRP.Range("L34").Value = 3827.389
$wsCRP.Range("M34").Value = -1753.0667
$wsCRP.Range("N34").Value = -4231.389
# Row 35
$wsCRP.Range("H35").Value = 62502420
$wsCRP.Range("I35").Value = 111111690
$wsCRP.Range("J35").Value = 4785.7144
$wsCRP.Range("K35").Value = 111111690
$wsCRP.Range("L35").Value = 4785.7144
$wsCRP.Range("M35").Value = -111111396
$wsCRP.Range("N35").Value = -5373.7144
# Row 58
$wsCRP.Range("H58").Value = 3685.3845
$wsCRP.Range("I58").Value = 4794.231
$wsCRP.Range("J58").Value = 2576.5386
$wsCRP.Range("K58").Value = 4794.231
$wsCRP.Range("L58").Value = 2576.5386
$wsCRP.Range("M58").Value = -4591.231
$wsCRP.Range("N58").Value = -2982.5386
# Row 105
$wsCRP.Range("H105").Value = 1098.697
$wsCRP.Range("I105").Value = 948.36
$wsCRP.Range("J105").Value = 1568.5
$wsCRP.Range("K105").Value = 948.36
$wsCRP.Range("L105").Value = 1568.5
$wsCRP.Range("M105").Value = 798.64
$wsCRP.Range("N105").Value = -5062.5
# Row 134
$wsCRP.Range("H134").Value = 1911.7241
$wsCRP.Range("I134").Value = 1212.3125
$wsCRP.Range("J134").Value = 2772.5386
$wsCRP.Range("K134").Value = 3636.9375
$wsCRP.Range("L134").Value = 8317.6158
$wsCRP.Range("M134").Value = -1101.9375
$wsCRP.Range("N134").Value = -13387.6158
# Row 136
$wsCRP.Range("H136").Value = 3685.3845
$wsCRP.Range("I136").Value = 4794.231
$wsCRP.Range("J136").Value = 2576.5386
$wsCRP.Range("K136").Value = 14382.693
$wsCRP.Range("L136").Value = 7729.6158
$wsCRP.Range("M136").Value = -11832.693
$wsCRP.Range("N136").Value = -12829.6158

$wsCUL = $wb.Worksheets.Item("CUL")
# Row 56
$wsCUL.Range("H56").Value = 4620.6665
$wsCUL.Range("I56").Value = 4620.6665
$wsCUL.Range("K56").Value = 4620.6665
$wsCUL.Range("M56").Value = -4090.6665
# Row 92
$wsCUL.Range("H92").Value = 26315968
$wsCUL.Range("I92").Value = 29411928
$wsCUL.Range("K92").Value = 88235784
$wsCUL.Range("M92").Value = -88234536

$wsGSM = $wb.Worksheets.Item("GSM")
# Row 70
$wsGSM.Range("H70").Value = 5604.767
$wsGSM.Range("I70").Value = 4004.6047
$wsGSM.Range("J70").Value = 7898.3335
$wsGSM.Range("K70").Value = 4004.6047
$wsGSM.Range("L70").Value = 7898.3335
$wsGSM.Range("M70").Value = -3734.6047
$wsGSM.Range("N70").Value = -8438.333500000001
# Row 73
$wsGSM.Range("H73").Value = 5604.767
$wsGSM.Range("I73").Value = 4004.6047
$wsGSM.Range("J73").Value = 7898.3335
$wsGSM.Range("K73").Value = 4004.6047
$wsGSM.Range("L73").Value = 7898.3335
$wsGSM.Range("M73").Value = -3068.6047
$wsGSM.Range("N73").Value = -9770.333500000001

$wsWVR = $wb.Worksheets.Item("WVR")
# Row 122
$wsWVR.Range("H122").Value = 55557324
$wsWVR.Range("I122").Value = 58825228
$wsWVR.Range("K122").Value = 176475684
$wsWVR.Range("M122").Value = -176473234
